$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.106.20'
$ws.Range('E2').Value = '  +2.53%  '
$ws.Range('D3').Value = '2.405.08'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''560.44'
$ws.Range('E5').Value = '  +2.19%  '
$ws.Range('D6').Value = '''138.19'
$ws.Range('E6').Value = '  +4.71%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''0.587'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').Value = '2.402.40'
$ws.Range('E9').Value = '  +2.94%  '
$ws.Range('E10').Value = '  +2.78%  '
$ws.Range('D11').Value = '''5.68'
$ws.Range('E11').Value = '  +3.01%  '
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  +3.13%  '
$ws.Range('D14').Value = '''25.80'
$ws.Range('E14').Value = '  +8.16%  '
$ws.Range('D15').Value = '2.831.32'
$ws.Range('E15').Value = '  +2.80%  '
$ws.Range('D16').Value = '62.022.85'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('E17').Value = '  +3.92%  '
$ws.Range('D18').Value = '2.399.15'
$ws.Range('E18').Value = '  +2.77%  '
$ws.Range('D19').Value = '''11.03'
$ws.Range('E19').Value = '  +3.50%  '
$ws.Range('D20').Value = '''343.66'
$ws.Range('E20').Value = '  +8.89%  '
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('D22').Value = '''6.88'
$ws.Range('E22').Value = '  +3.09%  '
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').Value = '''65.15'
$ws.Range('E24').Value = '  +1.33%  '
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '''8.37'
$ws.Range('E27').Value = '  +6.39%  '
$ws.Range('D28').Value = '''1.51'
$ws.Range('E28').Value = '  +10.95%  '
$ws.Range('D29').Value = '''1.38'
$ws.Range('E29').Value = '  +14.03%  '
$ws.Range('E30').Value = '  +3.91%  '
$ws.Range('D31').Value = '0.0₃0773'
$ws.Range('E31').Value = '  +4.50%  '
$ws.Range('D32').Value = '''6.37'
$ws.Range('E32').Value = '  +6.95%  '
$ws.Range('D33').Value = '''171.67'
$ws.Range('E33').Value = '  -1.03%  '
$ws.Range('E34').Value = '  +2.82%  '
$ws.Range('D35').Value = '''0.395'
$ws.Range('E35').Value = '  +3.54%  '
$ws.Range('E36').Value = '  +3.32%  '
$ws.Range('D37').Value = '''4.52'
$ws.Range('E37').Value = '  +10.05%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '''361.99'
$ws.Range('E38').Value = '  +10.66%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').Value = '''0.998'
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('E41').Value = '  +8.49%  '
$ws.Range('D42').Value = '''39.05'
$ws.Range('E42').Value = '  +2.41%  '
$ws.Range('D43').Value = '''143.88'
$ws.Range('E43').Value = '  +2.71%  '
$ws.Range('D44').Value = '''3.67'
$ws.Range('E44').Value = '  +5.23%  '
$ws.Range('D45').Value = '''20.62'
$ws.Range('E45').Value = '  +6.20%  '
$ws.Range('D46').Value = '''0.0965'
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('D47').Value = '''0.0518'
$ws.Range('E47').Value = '  +3.91%  '
$ws.Range('D48').Value = '''0.583'
$ws.Range('E48').Value = '  +3.47%  '
$ws.Range('E49').Value = '  +3.61%  '
$ws.Range('D50').Value = '''17.87'
$ws.Range('E50').Value = '  +4.70%  '
$ws.Range('E51').Value = '  -5.45%  '
